$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "Datos actualizados" timestamp (row 1, col A)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 14 de Junio de 2020 a las 16:23"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 2143631
$ws.Cells.Item(4, 3).Value = 1407
$ws.Cells.Item(4, 5).Value = 1171983
$ws.Cells.Item(4, 7).Value = 15
$ws.Cells.Item(4, 8).Value = 117542

# Row 5: Brasil
$ws.Cells.Item(5, 2).Value = 851321
$ws.Cells.Item(5, 3).Value = 525
$ws.Cells.Item(5, 5).Value = 371007
$ws.Cells.Item(5, 7).Value = 11
$ws.Cells.Item(5, 8).Value = 42802

# Row 7: India
$ws.Cells.Item(7, 2).Value = 324482
$ws.Cells.Item(7, 3).Value = 2856
$ws.Cells.Item(7, 4).Value = 164530
$ws.Cells.Item(7, 5).Value = 150705
$ws.Cells.Item(7, 7).Value = 48
$ws.Cells.Item(7, 8).Value = 9247

# Row 32: Emiratos Arabes Unidos
$ws.Cells.Item(32, 2).Value = 42294
$ws.Cells.Item(32, 3).Value = 304
$ws.Cells.Item(32, 4).Value = 27462
$ws.Cells.Item(32, 5).Value = 14543
$ws.Cells.Item(32, 7).Value = 1
$ws.Cells.Item(32, 8).Value = 289

# Row 57: Serbia
$ws.Cells.Item(57, 4).Value = 11511
$ws.Cells.Item(57, 5).Value = 545

# Row 60: Moldavia
$ws.Cells.Item(60, 2).Value = 11740
$ws.Cells.Item(60, 3).Value = 281
$ws.Cells.Item(60, 4).Value = 6623
$ws.Cells.Item(60, 5).Value = 4711
$ws.Cells.Item(60, 7).Value = 8
$ws.Cells.Item(60, 8).Value = 406

# Row 68: Noruega
$ws.Cells.Item(68, 2).Value = 8629
$ws.Cells.Item(68, 3).Value = 1
$ws.Cells.Item(68, 5).Value = 249

# Row 99: Cuba
$ws.Cells.Item(99, 2).Value = 2248
$ws.Cells.Item(99, 3).Value = 10
$ws.Cells.Item(99, 4).Value = 1948
$ws.Cells.Item(99, 5).Value = 216

# Rows 206/207: swap Islas Malvinas <-> Groenlandia (entire row data)
$ws.Cells.Item(206, 1).Value = "Groenlandia"
$ws.Cells.Item(207, 1).Value = "Islas Malvinas"

# Rows 210/211: swap Seychelles <-> Montserrat (entire row data)
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# Rows 213/214: swap Islas Virgenes Britanicas <-> Papua Nueva Guinea (entire row data)
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 8).Value = 0

$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 8).Value = 1
